# Fix on course info
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the Department value for the last course row: remove the stray
# trailing newline and instead enable wrap text on the cell.
$cell = $ws.Range("C4")
$cell.Value = "Engineering"
$cell.WrapText = $true

# Move the active selection to the edited cell.
$ws.Range("C4").Select()
